# Update main flow comments and some other code minor changes
#
# - Refresh the captured Start/End/Run time strings on the "Version" sheet.
# - Update the averaged/measured values on the "3_3" and "1_8" result sheets
#   and drop the now-unused "BT Idle" raw-data row from each.

$wb = $excel.ActiveWorkbook

# --- Sheet "Version": refresh timestamps -------------------------------
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "23:03:07 03-26-2017 Sunday Pacific Daylight Time"
$wsVersion.Range("C11").Value = "23:03:24 03-26-2017 Sunday Pacific Daylight Time"
$wsVersion.Range("C12").Value = "0:0:16 0 days"

# --- Sheet "3_3": update Deep Sleep stats, remove BT Idle row ----------
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.5659999999999999
$ws33.Range("D3").Value = 0.601
$ws33.Range("E3").Value = 0.543
$ws33.Range("F3").Value = 0.01
$ws33.Range("G3").Value = 100
$ws33.Range("H3").Value = "0.558333,0.56298,0.565107,0.558422,0.555181,0.563442,0.578375,0.579489,0.556846,0.570025,0.573463,0.554339,0.581058,0.56027,0.567423,0.550528,0.56743,0.556751,0.561112,0.564733,0.600929,0.555263,0.564264,0.558055,0.567423,0.590453,0.562688,0.56334,0.56455,0.570019,0.566126,0.562315,0.562043,0.595644,0.567518,0.558048,0.561941,0.567695,0.590079,0.563993,0.57207,0.567512,0.562505,0.593877,0.561676,0.554434,0.556744,0.56222,0.568721,0.564828,0.572715,0.565113,0.58832,0.562872,0.552111,0.543381,0.568347,0.561112,0.559441,0.589054,0.565202,0.552471,0.5651,0.562973,0.574659,0.563245,0.568272,0.557953,0.568259,0.556295,0.573633,0.566024,0.557681,0.556187,0.568449,0.564455,0.559815,0.565936,0.565759,0.565657,0.558415,0.56046,0.560745,0.564373,0.562416,0.554978,0.557498,0.561751,0.567613,0.56919,0.569747,0.566398,0.586364,0.564821,0.556839,0.568266,0.556472,0.561669,0.582648,0.561764"
$ws33.Rows.Item(4).Delete()

# --- Sheet "1_8": update Deep Sleep stats, remove BT Idle row ----------
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0.002
$ws18.Range("D3").Value = 0.006
$ws18.Range("E3").Value = -0.008999999999999999
$ws18.Range("F3").Value = 0.003
$ws18.Range("G3").Value = 100
$ws18.Range("H3").Value = "-0.001903,0.001125,0.000317,-0.001878,-0.002452,0.002259,-0.002994,-0.002446,-0.001348,-0.004425,-0.005233,0.004467,-0.006602,-0.003562,-0.006078,-0.00884,-0.001891,0.001717,0.00168,-0.004129,-0.002742,-0.002187,0.000866,-0.002994,-0.006361,0.000878,-0.00662,-0.002735,-0.000485,-0.008822,-0.003864,0.000329,-0.000787,0.001143,-0.003562,-0.001348,-0.004111,0.004195,-0.001903,-0.002476,-0.000528,-0.008285,0.002006,-0.005208,-0.00329,-0.002735,-0.0044,-0.000805,-0.000509,2e-06,-0.002187,-0.005239,-0.00136,-0.001323,0.001704,-0.002994,-0.00469,-0.003038,-0.00022,-0.00133,-0.005239,0.001168,0.000329,0.000588,-0.006602,-0.003297,-0.006596,0.001698,0.003369,3.9e-05,0.003647,-0.007181,0.00639,-0.002988,0.001427,0.003659,0.001686,-0.004407,-0.005227,-0.00329,-0.003833,0.003079,-0.003007,-0.003586,3.9e-05,-0.006053,-0.001891,0.000619,0.000594,-0.002735,-0.000799,-0.003852,-0.000775,-0.002187,-0.009383,-0.003013,-0.002181,-0.006922,0.00057,-0.001909"
$ws18.Rows.Item(4).Delete()
